# Update column G (K = strikeouts) values for rows 2-12
# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 2
    4  = 1
    5  = 2
    6  = 5
    7  = 2
    8  = 5
    9  = 8
    10 = 4
    11 = 1
    12 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
